# Generate Report for Handoff
# - Mark the handoff Priority as "ht" for the rows whose handback/priority
#   check failed (the "012a1e5f..." row and the later rows that share its
#   handoff batch) on both the zh-cn and de-de sheets.
# - Refresh the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#   timestamps for those same rows to reflect the newly generated report.

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14)

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-06 16:42:57"
}

# --- zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-06 16:42:52"
}

# --- de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-06 16:42:57"
}
